$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unhide columns E:AG (columns 5-33), which were hidden in a prior "session 3" cleanup
$ws.Range("E1:AG1").EntireColumn.Hidden = $false

# Update the active selection to F14
$ws.Range("F14").Select()
